# Updated 0 layer test cases
# - Fill in the missing Conf Matrix Acc / Loss Graph cells for run 166 (row 27)
# - Append a new test-case row for run 167 (row 29), mirroring row 27 but with
#   "Custom Preprocessing Used?" = No
# - Re-point the active selection the way the author left it (F14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 27: add the Conf Matrix Acc + Loss Graph columns that were pending ----
$ws.Range("L27").Value = 96.15
$ws.Range("M27").Value = "loss_graph.jpg"
$ws.Hyperlinks.Add($ws.Range("M27"), "..\graphs\train_metrics_166epochs_efficientnet_new\loss_graph.jpg") | Out-Null

# ---- Row 29: brand-new test case (run 167, 0 trainable layers, no custom preprocessing) ----
$ws.Range("A29").Value = 167
$ws.Range("B29").Value = "EfficientNet"
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = "No"
$ws.Range("E29").Value = "Shift, rotate, flip, shear, zoom, brightness"
$ws.Range("F29").Value = "RMSprop, epsilon = 1"
$ws.Range("G29").Value = 0.01
$ws.Range("H29").Value = "0.94 every 2 epochs"
$ws.Range("I29").Value = 16
$ws.Range("J29").Value = 0.01
$ws.Range("K29").Value = 0.2
$ws.Range("L29").Value = 96.62
$ws.Range("M29").Value = "loss_graph.jpg"
$ws.Hyperlinks.Add($ws.Range("M29"), "..\graphs\train_metrics_167epochs_efficientnet_new\loss_graph.jpg") | Out-Null

# Match column A's centered style used by every other row in the table.
$ws.Range("A27").Copy()
$ws.Range("A29").PasteSpecial(-4122)

# Match the rest of the row's left-aligned "Normal" style.
$ws.Range("B27:L27").Copy()
$ws.Range("B29:L29").PasteSpecial(-4122)

# Restore the Hyperlink cell style (the Hyperlinks.Add call above re-styles the
# cell with a plain underline font, not the workbook's named "Hyperlink" style)
# by pulling the format from an existing hyperlinked cell, without touching the
# hyperlink we just attached.
$ws.Range("M25").Copy()
$ws.Range("M27").PasteSpecial(-4122)
$ws.Range("M29").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Selection left where the author ended up ----
$ws.Range("F14").Select() | Out-Null
